$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=20.94432133333333; H=62.832964; I=0.7396577289668299; J=0.7396577289668298; K=3; M=216.1448186666667; N=648.434456; O=0.9739197284603751; P=0.9739197284603752; Q=4527.006536689732; R=40743.05883020758; S=0.7203672545489928; T=0.7203672545489926 }
    3  = @{ E=3; G=20.94432133333333; H=62.832964; I=0.7396577289668299; J=0.7396577289668298; K=3; M=3.181559666666666; N=9.544678999999999; O=0.01433568357434949; P=0.01433568357434949; Q=66.63560799983955; R=599.7204719985559; S=0.01060349915579043; T=0.01060349915579043 }
    4  = @{ E=3; G=20.94432133333333; H=62.832964; I=0.7396577289668299; J=0.7396577289668298; K=3; M=2.606510333333333; N=7.819531; O=0.01174458796527538; P=0.01174458796527538; Q=54.59158997998711; R=491.324309819884; S=0.008686975262046747; T=0.008686975262046747 }
    5  = @{ E=3; G=2.327094666666667; H=6.981284; I=0.08218235047311259; J=0.08218235047311258; K=3; M=216.1448186666667; N=648.434456; O=0.9739197284603751; P=0.9739197284603752; Q=502.9894547468338; R=4526.905092721504; S=0.0800390124570092; T=0.08003901245700919 }
    6  = @{ E=3; G=2.327094666666667; H=6.981284; I=0.08218235047311259; J=0.08218235047311258; K=3; M=3.181559666666666; N=9.544678999999999; O=0.01433568357434949; P=0.01433568357434949; Q=7.403790531981778; R=66.634114787836; S=0.001178140171778833; T=0.001178140171778833 }
    7  = @{ E=3; G=2.327094666666667; H=6.981284; I=0.08218235047311259; J=0.08218235047311258; K=3; M=2.606510333333333; N=7.819531; O=0.01174458796527538; P=0.01174458796527538; Q=6.065596295311556; R=54.590366657804; S=0.0009651978443245612; T=0.0009651978443245612 }
    8  = @{ E=3; G=5.044818; H=15.134454; I=0.1781599205600575; J=0.1781599205600575; K=3; M=216.1448186666667; N=648.434456; O=0.9739197284603751; P=0.9739197284603752; Q=1090.411271816336; R=9813.701446347024; S=0.1735134614543732; T=0.1735134614543732 }
    9  = @{ E=3; G=5.044818; H=15.134454; I=0.1781599205600575; J=0.1781599205600575; K=3; M=3.181559666666666; N=9.544678999999999; O=0.01433568357434949; P=0.01433568357434949; Q=16.050389474474; R=144.453505270266; S=0.002554044246780226; T=0.002554044246780226 }
    10 = @{ E=3; G=5.044818; H=15.134454; I=0.1781599205600575; J=0.1781599205600575; K=3; M=2.606510333333333; N=7.819531; O=0.01174458796527538; P=0.01174458796527538; Q=13.149370246786; R=118.344332221074; S=0.002092414858904069; T=0.002092414858904069 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
